# Pull & Bear crawling
# Add 5 new worksheets (one per store) for the "29 - 9" crawl run, each
# with the same Hora/Cambió/Nuevos/Actualizados layout used by the other
# sheets in this workbook, and the observed timestamps for that run.

$wb = $excel.ActiveWorkbook

$storeNames = @("29 - 9 PullAndBear", "29 - 9 Mango", "29 - 9 Zara", "29 - 9 Stradivarius", "29 - 9 Bershka")

$times = @("17:13", "17:16", "17:41", "17:52", "17:53", "17:55", "17:56", "18:8", "18:24", "18:27", "20:1", "20:8")

# Use the last existing sheet's header row as the formatting template
# (bold font, borders, centered alignment == style index 1).
$headerTemplate = $wb.Worksheets.Item($wb.Worksheets.Count)

foreach ($storeName in $storeNames) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $lastSheet)
    $ws.Name = $storeName

    $headerTemplate.Range("A1:D1").Copy($ws.Range("A1:D1"))
    $ws.Range("A1").Value = "Hora"
    $ws.Range("B1").Value = "Cambió"
    $ws.Range("C1").Value = "Nuevos"
    $ws.Range("D1").Value = "Actualizados"

    $row = 2
    foreach ($t in $times) {
        $ws.Cells.Item($row, 1).Value = $t
        $ws.Cells.Item($row, 2).Value = $false
        $ws.Cells.Item($row, 3).Value = 0
        $ws.Cells.Item($row, 4).Value = 0
        $row = $row + 1
    }
}
